# refactor: merge prefabs for login and signup into one
# Also created a prefab for the promo code
#
# Adds a new "UI Promo code" localization row (row 29) to the
# tsv_UI_Defaults sheet, mirroring the existing "UI X" header-row pattern
# (Notes | Key | English | Français).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (28) onto the new
# row (29) so the new row picks up the same style (borders/wrap/valign)
# as every other data row, then overwrite the copied values.
$ws.Range("A28:F28").Copy()
$ws.Range("A29:F29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A29").Value = "UI Promo code"
$ws.Range("B29").Value = "tmp promo code"
$ws.Range("C29").Value = "Promo Code"
$ws.Range("D29").Value = "Code Promotionnel"

# Match the author's final selection/scroll position captured in the diff.
$ws.Range("D31").Select()
